$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96 (G96=19894)
$ws.Range("H96").Value = 1700
$ws.Range("I96").Value = 400
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 9000
$ws.Range("M96").Value = 173
$ws.Range("N96").Value = -11746

# Row 125 (G125=36228)
$ws.Range("H125").Value = 2739.4375
$ws.Range("I125").Value = 1599.4286
$ws.Range("J125").Value = 3626.111
$ws.Range("K125").Value = 14394.8574
$ws.Range("L125").Value = 32634.999
$ws.Range("M125").Value = -11934.8574
$ws.Range("N125").Value = -37554.999

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (G45=27714)
$ws.Range("H45").Value = 1735
$ws.Range("I45").Value = 1091.4814
$ws.Range("J45").Value = 5210
$ws.Range("K45").Value = 1091.4814
$ws.Range("L45").Value = 5210
$ws.Range("M45").Value = -714.4813999999999
$ws.Range("N45").Value = -5964

# Row 122 (G122=36168)
$ws.Range("H122").Value = 2825.1538
$ws.Range("I122").Value = 1968.375
$ws.Range("J122").Value = 4196
$ws.Range("K122").Value = 5905.125
$ws.Range("L122").Value = 12588
$ws.Range("M122").Value = -3455.125
$ws.Range("N122").Value = -17488

# Row 123 (G123=34107)
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (G94=19939)
$ws.Range("H94").Value = 618.12
$ws.Range("I94").Value = 598.2381
$ws.Range("J94").Value = 722.5
$ws.Range("K94").Value = 598.2381
$ws.Range("L94").Value = 722.5
$ws.Range("M94").Value = -147.2381
$ws.Range("N94").Value = -1624.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G31=44023)
$ws.Range("H31").Value = 2002433.9
$ws.Range("I31").Value = 2274379.5
$ws.Range("J31").Value = 8166.6665
$ws.Range("K31").Value = 2274379.5
$ws.Range("L31").Value = 8166.6665
$ws.Range("M31").Value = -2274084.5
$ws.Range("N31").Value = -8756.666499999999

# Row 34 (G34=44023)
$ws.Range("H34").Value = 2002433.9
$ws.Range("I34").Value = 2274379.5
$ws.Range("J34").Value = 8166.6665
$ws.Range("K34").Value = 2274379.5
$ws.Range("L34").Value = 8166.6665
$ws.Range("M34").Value = -2274177.5
$ws.Range("N34").Value = -8570.666499999999

# Row 41 (G41=1917)
$ws.Range("H41").Value = 4269.3
$ws.Range("I41").Value = 2318.6
$ws.Range("J41").Value = 6220
$ws.Range("K41").Value = 2318.6
$ws.Range("L41").Value = 6220
$ws.Range("M41").Value = -1890.6
$ws.Range("N41").Value = -7076

# Row 50 (G50=1862)
$ws.Range("H50").Value = 14400
$ws.Range("J50").Value = 14400
$ws.Range("L50").Value = 14400
$ws.Range("N50").Value = -15650

# Row 51 (G51=2039)
$ws.Range("H51").Value = 13296.667
$ws.Range("I51").Value = 1980
$ws.Range("J51").Value = 15560
$ws.Range("K51").Value = 1980
$ws.Range("L51").Value = 15560
$ws.Range("M51").Value = -1244
$ws.Range("N51").Value = -17032

# Row 58 (G58=44021)
$ws.Range("H58").Value = 17243844
$ws.Range("I58").Value = 1045.7273
$ws.Range("J58").Value = 27781110
$ws.Range("K58").Value = 1045.7273
$ws.Range("L58").Value = 27781110
$ws.Range("M58").Value = -842.7273
$ws.Range("N58").Value = -27781516

# Row 59 (G59=1942)
$ws.Range("H59").Value = 17733.334
$ws.Range("J59").Value = 17733.334
$ws.Range("L59").Value = 17733.334
$ws.Range("N59").Value = -20023.334

# Row 60 (G60=1937)
$ws.Range("H60").Value = 11677.08
$ws.Range("I60").Value = 16000
$ws.Range("J60").Value = 11496.958
$ws.Range("K60").Value = 16000
$ws.Range("L60").Value = 11496.958
$ws.Range("M60").Value = -15489
$ws.Range("N60").Value = -12518.958

# Row 61 (G61=2039)
$ws.Range("H61").Value = 13296.667
$ws.Range("I61").Value = 1980
$ws.Range("J61").Value = 15560
$ws.Range("K61").Value = 1980
$ws.Range("L61").Value = 15560
$ws.Range("M61").Value = -1632
$ws.Range("N61").Value = -16256

# Row 62 (G62=12580)
$ws.Range("H62").Value = 5231.1113
$ws.Range("I62").Value = 2860
$ws.Range("J62").Value = 6416.6665
$ws.Range("K62").Value = 2860
$ws.Range("L62").Value = 6416.6665
$ws.Range("M62").Value = -2236
$ws.Range("N62").Value = -7664.6665

# Row 65 (G65=12580)
$ws.Range("H65").Value = 5231.1113
$ws.Range("I65").Value = 2860
$ws.Range("J65").Value = 6416.6665
$ws.Range("K65").Value = 14300
$ws.Range("L65").Value = 32083.3325
$ws.Range("M65").Value = -11180
$ws.Range("N65").Value = -38323.3325

# Row 68 (G68=10611)
$ws.Range("H68").Value = 29500
$ws.Range("J68").Value = 29500
$ws.Range("L68").Value = 29500
$ws.Range("N68").Value = -30998

# Row 71 (G71=10611)
$ws.Range("H71").Value = 29500
$ws.Range("J71").Value = 29500
$ws.Range("L71").Value = 88500
$ws.Range("N71").Value = -95988

# Row 74 (G74=10636)
$ws.Range("H74").Value = 23299.445
$ws.Range("J74").Value = 23299.445
$ws.Range("L74").Value = 23299.445
$ws.Range("N74").Value = -25047.445

# Row 77 (G77=10636)
$ws.Range("H77").Value = 23299.445
$ws.Range("J77").Value = 23299.445
$ws.Range("L77").Value = 69898.33499999999
$ws.Range("N77").Value = -78634.33499999999

# Row 136 (G136=44021)
$ws.Range("H136").Value = 17243844
$ws.Range("I136").Value = 1045.7273
$ws.Range("J136").Value = 27781110
$ws.Range("K136").Value = 3137.1819
$ws.Range("L136").Value = 83343330
$ws.Range("M136").Value = -587.1819
$ws.Range("N136").Value = -83348430

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (G5=43974)
$ws.Range("H5").Value = 1294.9166
$ws.Range("I5").Value = 498
$ws.Range("J5").Value = 2888.75
$ws.Range("K5").Value = 1494
$ws.Range("L5").Value = 8666.25
$ws.Range("M5").Value = -1382
$ws.Range("N5").Value = -8890.25

# Row 135 (G135=43974)
$ws.Range("H135").Value = 1294.9166
$ws.Range("I135").Value = 498
$ws.Range("J135").Value = 2888.75
$ws.Range("K135").Value = 4482
$ws.Range("L135").Value = 25998.75
$ws.Range("M135").Value = -1947
$ws.Range("N135").Value = -31068.75

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (G122=36247)
$ws.Range("H122").Value = 2829.244
$ws.Range("I122").Value = 2525.3215
$ws.Range("K122").Value = 7575.9645
$ws.Range("M122").Value = -5125.9645

$ws = $wb.Worksheets.Item("WVR")
# Row 16 (G16=26304)
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 107 (G107=27746)
$ws.Range("H107").Value = 1555
$ws.Range("I107").Value = 1109.0625
$ws.Range("J107").Value = 3933.3333
$ws.Range("K107").Value = 3327.1875
$ws.Range("L107").Value = 11799.9999
$ws.Range("M107").Value = -1407.1875
$ws.Range("N107").Value = -15639.9999

# Row 132 (G132=44029)
$ws.Range("H132").Value = 598217.9
$ws.Range("I132").Value = 1253713.2
$ws.Range("J132").Value = 15555.333
$ws.Range("K132").Value = 3761139.6
$ws.Range("L132").Value = 46665.999
$ws.Range("M132").Value = -3758609.6
$ws.Range("N132").Value = -51725.999

# Row 136 (G136=44031)
$ws.Range("H136").Value = 1905.3334
$ws.Range("I136").Value = 983
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 2949
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -399
$ws.Range("N136").Value = -16350

# Row 141 (G141=42505)
$ws.Range("H141").Value = 26660
$ws.Range("J141").Value = 26660
$ws.Range("L141").Value = 26660
$ws.Range("N141").Value = -37020
